$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "mora" period values: the 4 periods (2507,2506,2505,2504) become
# 3 periods (2506,2507,2508) - the whole window shifted forward by one period.
$ws.Range("E16").Value = "2506"
$ws.Range("E17").Value = "2507"

# Row 18 (period 2505) is removed entirely; this shifts the old row 19
# (period 2504, with its own distinct border/style) up to become the new
# row 18, carrying its own formatting with it.
$ws.Rows("18").Delete()

# The row that is now 18 (previously row 19) needs its period updated to
# the new latest period, 2508.
$ws.Range("E18").Value = "2508"

# Cant. Periodos (count of periods) drops from 4 to 3.
$ws.Range("F13").Value = 3

# Valor Mora total updates accordingly (56940 * 3 = 170820).
$ws.Range("E11").Value = 170820
